$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Moldova now reports 1504 patients; its highlight becomes green (completed)
$ws.Range("B9").Value = 1504
$ws.Range("B9").Interior.Color = 5287936

# Kazakhstan is now marked as in-progress (orange highlight, no bold)
$ws.Range("B7").Interior.Color = 49407

# Update the active selection
$ws.Range("G13").Select()
